# ------------------------------------------------------------------------
# Append scrape snapshot: 2025-09-16 12:36 JST run
#
# The scraper re-ran and produced a refreshed snapshot of the listings.
# New items were prepended/inserted, one stale item was dropped, and new
# items were appended, so every data row shifts. Rather than replaying
# that insert/delete logic cell by cell, we rewrite the whole data block
# (rows 2-15) with the final known values, refresh the timestamp column,
# widen column H, and rebuild the hyperlinks on column F.
# ------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Widen column H (8th column) from stored width 12 to stored width 18 ---
# Excel's COM ColumnWidth property is expressed in "characters" and is
# offset from the raw OOXML stored width by a constant +5/6 (default
# Calibri 11 padding); 17.166666666666668 round-trips to a saved width of
# exactly 18.
$ws.Columns.Item(8).ColumnWidth = 17.166666666666668

# --- Drop all existing hyperlinks on the sheet; they will be rebuilt from
#     scratch below once every F-column cell has its final URL. ---
$ws.Range("A1").Hyperlinks.Delete()

# --- Rewrite data rows 2-15 with the refreshed scrape contents ---

# Row 2: 【急募】AI&SaaS Lineプラットフォーム開発のプロを探しています!
$ws.Cells.Item(2,1).Value = '2025-09-16 12:36:17'
$ws.Cells.Item(2,2).Value = '【急募】AI&SaaS Lineプラットフォーム開発のプロを探しています!'
$ws.Cells.Item(2,3).Value = 'システム開発'
$ws.Cells.Item(2,4).Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Cells.Item(2,5).Value = '期限情報なし'
$ws.Cells.Item(2,7).Value = 375
$ws.Cells.Item(2,8).Value = '🔥AI,Ai ◆開発'
$ws.Cells.Item(2,6).Value = 'https://www.lancers.jp/work/detail/5393834'
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5393834') | Out-Null
$ws.Cells.Item(2,6).Style = "Hyperlink"

# Row 3: 【AI活用】データ分析Webサービス開発パートナー募集
$ws.Cells.Item(3,1).Value = '2025-09-16 12:36:17'
$ws.Cells.Item(3,2).Value = '【AI活用】データ分析Webサービス開発パートナー募集'
$ws.Cells.Item(3,3).Value = 'システム開発'
$ws.Cells.Item(3,4).Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Cells.Item(3,5).Value = '期限情報なし'
$ws.Cells.Item(3,7).Value = 368
$ws.Cells.Item(3,8).Value = '🔥AI,Ai ◆開発'
$ws.Cells.Item(3,6).Value = 'https://www.lancers.jp/work/detail/5393929'
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5393929') | Out-Null
$ws.Cells.Item(3,6).Style = "Hyperlink"

# Row 4: 【急募】AIを活用した製造・輸入計画最適化システム構築
$ws.Cells.Item(4,1).Value = '2025-09-16 12:36:17'
$ws.Cells.Item(4,2).Value = '【急募】AIを活用した製造・輸入計画最適化システム構築'
$ws.Cells.Item(4,3).Value = 'システム開発'
$ws.Cells.Item(4,4).Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(4,5).Value = '期限情報なし'
$ws.Cells.Item(4,7).Value = 318
$ws.Cells.Item(4,8).Value = '🔥AI,Ai'
$ws.Cells.Item(4,6).Value = 'https://www.lancers.jp/work/detail/5394475'
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5394475') | Out-Null
$ws.Cells.Item(4,6).Style = "Hyperlink"

# Row 5: AIチャットボットのβ版テスト参加者募集!
$ws.Cells.Item(5,1).Value = '2025-09-16 12:36:17'
$ws.Cells.Item(5,2).Value = 'AIチャットボットのβ版テスト参加者募集!'
$ws.Cells.Item(5,3).Value = 'システム開発'
$ws.Cells.Item(5,4).Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Cells.Item(5,5).Value = '期限情報なし'
$ws.Cells.Item(5,7).Value = 295
$ws.Cells.Item(5,8).Value = '🔥AI,Ai'
$ws.Cells.Item(5,6).Value = 'https://www.lancers.jp/work/detail/5394484'
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5394484') | Out-Null
$ws.Cells.Item(5,6).Style = "Hyperlink"

# Row 6: 【せどり×ツール製作】APIを使用したせどりツールを製作できるエンジニアさんを募
$ws.Cells.Item(6,1).Value = '2025-09-16 12:36:17'
$ws.Cells.Item(6,2).Value = '【せどり×ツール製作】APIを使用したせどりツールを製作できるエンジニアさんを募集します♪'
$ws.Cells.Item(6,3).Value = 'システム開発'
$ws.Cells.Item(6,4).Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(6,5).Value = '期限情報なし'
$ws.Cells.Item(6,7).Value = 243
$ws.Cells.Item(6,8).Value = '🔥API ◆ツール'
$ws.Cells.Item(6,6).Value = 'https://www.lancers.jp/work/detail/5217096'
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5217096') | Out-Null
$ws.Cells.Item(6,6).Style = "Hyperlink"

# Row 7: 【急募】大規模データ収集自動化(スクレイピング・DB連携・エラー管理)案件
$ws.Cells.Item(7,1).Value = '2025-09-16 12:36:17'
$ws.Cells.Item(7,2).Value = '【急募】大規模データ収集自動化(スクレイピング・DB連携・エラー管理)案件'
$ws.Cells.Item(7,3).Value = 'システム開発'
$ws.Cells.Item(7,4).Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Cells.Item(7,5).Value = '期限情報なし'
$ws.Cells.Item(7,7).Value = 158
$ws.Cells.Item(7,8).Value = '◆自動化,スクレイピング ◇管理'
$ws.Cells.Item(7,6).Value = 'https://www.lancers.jp/work/detail/5394578'
$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5394578') | Out-Null
$ws.Cells.Item(7,6).Style = "Hyperlink"

# Row 8: 【新規開店】新店舗自動抽出ツールの開発依頼
$ws.Cells.Item(8,1).Value = '2025-09-16 12:36:17'
$ws.Cells.Item(8,2).Value = '【新規開店】新店舗自動抽出ツールの開発依頼'
$ws.Cells.Item(8,3).Value = 'システム開発'
$ws.Cells.Item(8,4).Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Cells.Item(8,5).Value = '期限情報なし'
$ws.Cells.Item(8,7).Value = 128
$ws.Cells.Item(8,8).Value = '◆ツール,開発'
$ws.Cells.Item(8,6).Value = 'https://www.lancers.jp/work/detail/5394572'
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5394572') | Out-Null
$ws.Cells.Item(8,6).Style = "Hyperlink"

# Row 9: 【急募】微生物関連データ管理システムのサポート依頼
$ws.Cells.Item(9,1).Value = '2025-09-16 12:36:17'
$ws.Cells.Item(9,2).Value = '【急募】微生物関連データ管理システムのサポート依頼'
$ws.Cells.Item(9,3).Value = 'システム開発'
$ws.Cells.Item(9,4).Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Cells.Item(9,5).Value = '期限情報なし'
$ws.Cells.Item(9,7).Value = 45
$ws.Cells.Item(9,8).Value = '◇管理'
$ws.Cells.Item(9,6).Value = 'https://www.lancers.jp/work/detail/5394505'
$ws.Hyperlinks.Add($ws.Range("F9"), 'https://www.lancers.jp/work/detail/5394505') | Out-Null
$ws.Cells.Item(9,6).Style = "Hyperlink"

# Row 10: 【急募】モバイルアプリ・webアプリのバグ修正と機能実装
$ws.Cells.Item(10,1).Value = '2025-09-16 12:36:17'
$ws.Cells.Item(10,2).Value = '【急募】モバイルアプリ・webアプリのバグ修正と機能実装'
$ws.Cells.Item(10,3).Value = 'システム開発'
$ws.Cells.Item(10,4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(10,5).Value = '期限情報なし'
$ws.Cells.Item(10,7).Value = 45
$ws.Cells.Item(10,8).Value = '◇アプリ'
$ws.Cells.Item(10,6).Value = 'https://www.lancers.jp/work/detail/5394097'
$ws.Hyperlinks.Add($ws.Range("F10"), 'https://www.lancers.jp/work/detail/5394097') | Out-Null
$ws.Cells.Item(10,6).Style = "Hyperlink"

# Row 11: 【急募】Firebaseを活用したモバイル・Webアプリのログイン機能実装
$ws.Cells.Item(11,1).Value = '2025-09-16 12:36:17'
$ws.Cells.Item(11,2).Value = '【急募】Firebaseを活用したモバイル・Webアプリのログイン機能実装'
$ws.Cells.Item(11,3).Value = 'システム開発'
$ws.Cells.Item(11,4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(11,5).Value = '期限情報なし'
$ws.Cells.Item(11,7).Value = 45
$ws.Cells.Item(11,8).Value = '◇アプリ'
$ws.Cells.Item(11,6).Value = 'https://www.lancers.jp/work/detail/5394060'
$ws.Hyperlinks.Add($ws.Range("F11"), 'https://www.lancers.jp/work/detail/5394060') | Out-Null
$ws.Cells.Item(11,6).Style = "Hyperlink"

# Row 12: 既存LSOサイト(ver.6)の更新
$ws.Cells.Item(12,1).Value = '2025-09-16 12:36:17'
$ws.Cells.Item(12,2).Value = '既存LSOサイト(ver.6)の更新'
$ws.Cells.Item(12,3).Value = 'システム開発'
$ws.Cells.Item(12,4).Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Cells.Item(12,5).Value = '期限情報なし'
$ws.Cells.Item(12,7).Value = 38
$ws.Cells.Item(12,8).Value = '◇サイト'
$ws.Cells.Item(12,6).Value = 'https://www.lancers.jp/work/detail/5394061'
$ws.Hyperlinks.Add($ws.Range("F12"), 'https://www.lancers.jp/work/detail/5394061') | Out-Null
$ws.Cells.Item(12,6).Style = "Hyperlink"

# Row 13: 自社のWebデータベースExmentの表示速度アップのご依頼
$ws.Cells.Item(13,1).Value = '2025-09-16 12:36:17'
$ws.Cells.Item(13,2).Value = '自社のWebデータベースExmentの表示速度アップのご依頼'
$ws.Cells.Item(13,3).Value = 'システム開発'
$ws.Cells.Item(13,4).Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Cells.Item(13,5).Value = '期限情報なし'
$ws.Cells.Item(13,7).Value = 25
$ws.Cells.Item(13,6).Value = 'https://www.lancers.jp/work/detail/5394186'
$ws.Hyperlinks.Add($ws.Range("F13"), 'https://www.lancers.jp/work/detail/5394186') | Out-Null
$ws.Cells.Item(13,6).Style = "Hyperlink"

# Row 14: 【自走型組織】サイボウズで売上・数量・粗利・経費・在庫を自動集計しグラフ化
$ws.Cells.Item(14,1).Value = '2025-09-16 12:36:17'
$ws.Cells.Item(14,2).Value = '【自走型組織】サイボウズで売上・数量・粗利・経費・在庫を自動集計しグラフ化'
$ws.Cells.Item(14,3).Value = 'システム開発'
$ws.Cells.Item(14,4).Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(14,5).Value = '期限情報なし'
$ws.Cells.Item(14,7).Value = 18
$ws.Cells.Item(14,6).Value = 'https://www.lancers.jp/work/detail/5394424'
$ws.Hyperlinks.Add($ws.Range("F14"), 'https://www.lancers.jp/work/detail/5394424') | Out-Null
$ws.Cells.Item(14,6).Style = "Hyperlink"

# Row 15: 【VBA】エクセルでのマクロ作成【excel】
$ws.Cells.Item(15,1).Value = '2025-09-16 12:36:17'
$ws.Cells.Item(15,2).Value = '【VBA】エクセルでのマクロ作成【excel】'
$ws.Cells.Item(15,3).Value = 'システム開発'
$ws.Cells.Item(15,4).Value = '~ 5,000 円 / 固定'
$ws.Cells.Item(15,5).Value = '期限情報なし'
$ws.Cells.Item(15,7).Value = 10
$ws.Cells.Item(15,6).Value = 'https://www.lancers.jp/work/detail/5394416'
$ws.Hyperlinks.Add($ws.Range("F15"), 'https://www.lancers.jp/work/detail/5394416') | Out-Null
$ws.Cells.Item(15,6).Style = "Hyperlink"
